$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the product title in cell B3:
# "هاهای" -> "های" ("کرم هاهای" -> "کرم های")
$ws.Range("B3").Value2 = "تولید کرم های گیاهی نرم کننده و مرطوب کننده پوست و .... و پماد ضد درد گیاهی"

# Move the active selection, matching the saved cursor position after the edit
$ws.Range("B8").Select()
